$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refreshed coin price snapshot (column D, "Price") pulled in by the
# scheduled GitHub Actions symbol-list update. Values are written back
# as text (NumberFormat "@") so the exact textual representation -
# including trailing zeros / significant digits - is preserved just
# like the original inline-string cells.
$priceUpdates = @{
    "D2"  = "248.75"
    "D3"  = "22.46"
    "D4"  = "5.395"
    "D5"  = "0.05693"
    "D6"  = "3.401"
    "D7"  = "6.319"
    "D8"  = "0.8067"
    "D9"  = "0.9132"
    "D11" = "0.07415"
    "D12" = "0.03129"
    "D13" = "0.03012"
    "D14" = "0.09380"
    "D15" = "3.859"
    "D16" = "0.001571"
    "D17" = "0.04769"
    "D18" = "0.01826"
    "D19" = "0.0005846"
    "D20" = "0.006450"
    "D21" = "0.004996"
    "D22" = "0.001005"
    "D23" = "0.0001499"
    "D24" = "3.702"
    "D25" = "2.199"
    "D26" = "0.3256"
    "D27" = "0.1306"
    "D40" = "0.04006"
    "D41" = "0.006828"
    "D43" = "0.002709"
    "D44" = "0.007514"
    "D45" = "0.00005757"
    "D46" = "0.00000000749"
    "D47" = "0.4986"
    "D48" = "0.2090"
    "D49" = "0.00002098"
    "D50" = "0.01009"
}

foreach ($cellRef in $priceUpdates.Keys) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$cellRef]
}
